$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2019-12-31 00:00:00"
$ws.Range("O2").Value = 139608187.69
$ws.Range("P2").Value = 1323.1237634968
$ws.Range("Q2").Value = 691610981.8099999
$ws.Range("R2").Value = 6554.6794945874
$ws.Range("S2").Value = 46252355.07
$ws.Range("T2").Value = 438.3524428144
$ws.Range("U2").Value = -84728304.8
$ws.Range("V2").Value = -803.0047189683
$ws.Range("W2").ClearContents()
$ws.Range("X2").ClearContents()
$ws.Range("Y2").Value = 85102852.41
$ws.Range("Z2").Value = 806.5544595068
$ws.Range("AA2").Value = -47463516.61
$ws.Range("AB2").Value = -449.8311149577
$ws.Range("AC2").Value = 10551408.08
$ws.Range("AD2").Value = 579.1432322156001
